# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder country names in the "Corea del Sur / Filipinas / Colombia" block ---
# Before: A41=Corea del Sur, A42=Filipinas, A43=Colombia
# After : A41=Colombia,      A42=Corea del Sur, A43=Filipinas
$ws.Range("A41").Value = "Colombia"
$ws.Range("A42").Value = "Corea del Sur"
$ws.Range("A43").Value = "Filipinas"

# --- Reorder country names in the "Republica de Yibuti / Sudan" block ---
# Before: A91=Republica de Yibuti, A92=Sudan
# After : A91=Sudan,               A92=Republica de Yibuti
$ws.Range("A91").Value = "Sudan"
$ws.Range("A92").Value = "Republica de Yibuti"

# --- Update data values (row 4 - Estados Unidos) ---
$ws.Range("B4").Value = 1366753
$ws.Range("C4").Value = 19444
$ws.Range("E4").Value = 1030342
$ws.Range("G4").Value = 720
$ws.Range("H4").Value = 80757

# --- Update data values (row 41 - now Colombia, new data) ---
$ws.Range("B41").Value = 11063
$ws.Range("C41").Value = 568
$ws.Range("D41").Value = 2705
$ws.Range("E41").Value = 7895
$ws.Range("F41").Value = 130
$ws.Range("G41").Value = 18
$ws.Range("H41").Value = 463

# --- Update data values (row 42 - now Corea del Sur, shifted from old row41 data) ---
$ws.Range("B42").Value = 10874
$ws.Range("C42").Value = 34
$ws.Range("D42").Value = 9610
$ws.Range("E42").Value = 1008
$ws.Range("F42").Value = 55
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 256

# --- Update data values (row 43 - now Filipinas, shifted from old row42 data) ---
$ws.Range("B43").Value = 10794
$ws.Range("C43").Value = 184
$ws.Range("D43").Value = 1924
$ws.Range("E43").Value = 8151
$ws.Range("F43").Value = 31
$ws.Range("G43").Value = 15
$ws.Range("H43").Value = 719

# --- Update data values (row 91 - now Sudan, new data) ---
$ws.Range("B91").Value = 1365
$ws.Range("C91").Value = 201
$ws.Range("D91").Value = 149
$ws.Range("E91").Value = 1146
$ws.Range("G91").Value = 6
$ws.Range("H91").Value = 70

# --- Update data values (row 92 - now Republica de Yibuti, shifted from old row91 data) ---
$ws.Range("B92").Value = 1210
$ws.Range("C92").Value = 21
$ws.Range("D92").Value = 847
$ws.Range("E92").Value = 360
$ws.Range("H92").Value = 3

# --- Update data values (row 113 - Crucero) ---
$ws.Range("D113").Value = 651
$ws.Range("E113").Value = 48
